$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet2 "burndown" data table: add Sprint-3 (col I) and Sprint-4 (col J)
# day-remaining figures for rows 96-124, then refresh the SUM totals in
# row 127 that the Sprint 4 Burndown chart reads from.
# ---------------------------------------------------------------------------
$newData = @{
    96  = @(2, 0)
    97  = @(0, 0)
    98  = @(1, 0)
    99  = @(4, 0)
    100 = @(4, 0)
    102 = @(0, 0)
    103 = @(0, 0)
    104 = @(0, 0)
    105 = @(0, 0)
    106 = @(0, 0)
    107 = @(0, 0)
    108 = @(0, 0)
    109 = @(0, 0)
    110 = @(3, 0)
    111 = @(4, 4)
    115 = @(0, 0)
    116 = @(0.5, 0)
    117 = @(1, 0)
    119 = @(4, 0)
    120 = @(4, 6)
    121 = @(1, 1)
    122 = @(1, 1)
    123 = @(0, 0)
    124 = @(0, 0)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws2.Cells.Item($row, 9).Value = $vals[0]
    $ws2.Cells.Item($row, 10).Value = $vals[1]
}

$ws2.Range("I127").Formula = "=SUM(I96:I124)"
$ws2.Range("J127").Formula = "=SUM(J96:J124)"

# ---------------------------------------------------------------------------
# Chart 5 ("Sprint 4 Burndown") on Sheet2: give it a real chart title,
# extend its category/value series references from G:H to G:J so the two
# new sprint columns show up, and move/resize it to its new anchor.
# ---------------------------------------------------------------------------
$co = $ws2.ChartObjects(5)
$chart = $co.Chart

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Sprint 4 Burndown"

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(`"Sprint 4 Burndown`",Sheet2!`$G`$94:`$J`$94,Sheet2!`$G`$127:`$J`$127,1)"

$co.Left = 820.7919134473425
$co.Top = 1561.5000787401575
$co.Width = 320.625
$co.Height = 218.33330708661424

# ---------------------------------------------------------------------------
# View-state tweaks captured in the diff: where each sheet is scrolled to
# and which cell(s) are selected on each.
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 11
$ws1.Range("A34:XFD34").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 93
$excel.ActiveWindow.ScrollColumn = 5
$ws2.Range("L98").Select()
